# Updated cryptos list on Sat Aug 17 04:42:04 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-Cell "D2" "59.035.60"
Set-Cell "E2" "  +1.86%  "

# Row 3 - Ethereum
Set-Cell "D3" "2.588.53"
Set-Cell "E3" "  +0.67%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  -0.02%  "

# Row 5 - BNB
Set-Cell "D5" "522.55"

# Row 6 - Solana
Set-Cell "D6" "139.19"
Set-Cell "E6" "  -2.17%  "

# Row 7 - USDC
Set-Cell "E7" "  +0.02%  "

# Row 8 - XRP
Set-Cell "E8" "  +0.20%  "

# Row 9 - LidoStakedEther
Set-Cell "D9" "2.599.87"
Set-Cell "E9" "  +0.50%  "

# Row 10 - Toncoin
Set-Cell "E10" "  -2.66%  "

# Row 11 - Dogecoin
Set-Cell "E11" "  -0.15%  "

# Row 12 - Cardano
Set-Cell "E12" "  +1.42%  "

# Row 13 - TRON
Set-Cell "E13" "  +3.30%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-Cell "D14" "3.041.80"
Set-Cell "E14" "  +0.49%  "

# Row 15 - WrappedBTC
Set-Cell "D15" "58.975.14"
Set-Cell "E15" "  +1.81%  "

# Row 16 - Avalanche
Set-Cell "E16" "  +0.72%  "

# Row 17 - now ShibaInu (was WrappedEther)
Set-Cell "B17" "ShibaInu"
Set-Cell "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-Cell "D17" "0.0000133"
Set-Cell "E17" "  -0.53%  "

# Row 18 - now WrappedEther (was ShibaInu)
Set-Cell "B18" "WrappedEther"
Set-Cell "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-Cell "D18" "2.570.57"
Set-Cell "E18" "  +0.23%  "

# Row 19 - BitcoinCash
Set-Cell "D19" "338.77"
Set-Cell "E19" "  -0.65%  "

# Row 20 - Polkadot
Set-Cell "E20" "  +0.02%  "

# Row 21 - Chainlink
Set-Cell "D21" "10.10"
Set-Cell "E21" "  -1.22%  "

# Row 22 - Uniswap
Set-Cell "D22" "6.52"
Set-Cell "E22" "  +2.90%  "

# Row 23 - Dai
Set-Cell "E23" "  +0.06%  "

# Row 24 - Litecoin
Set-Cell "D24" "65.99"
Set-Cell "E24" "  +1.01%  "

# Row 25 - Kaspa
Set-Cell "E25" "  +1.27%  "

# Row 26 - Polygon
Set-Cell "E26" "  +0.58%  "

# Row 27 - Binance-PegBSC-USD
Set-Cell "E27" "  +0.09%  "

# Row 28 - InternetComputer(DFINITY)
Set-Cell "E28" "  +0.44%  "

# Row 29 - USDe
Set-Cell "E29" "  -0.01%  "

# Row 30 - PEPE
Set-Cell "D30" "0.0₃0725"
Set-Cell "E30" "  -2.51%  "

# Row 31 - Aptos
Set-Cell "D31" "5.92"
Set-Cell "E31" "  -4.70%  "

# Row 32 - PancakeSwap
Set-Cell "E32" "  +0.44%  "

# Row 33 - EthereumClassic
Set-Cell "D33" "18.69"
Set-Cell "E33" "  +0.09%  "

# Row 34 - Monero
Set-Cell "D34" "148.98"
Set-Cell "E34" "  -0.58%  "

# Row 35 - NEARProtocol
Set-Cell "E35" "  -0.33%  "

# Row 36 - ImmutableX
Set-Cell "E36" "  -1.27%  "

# Row 37 - OKB
Set-Cell "E37" "  +2.32%  "

# Row 38 - Stacks
Set-Cell "E38" "  +1.85%  "

# Row 39 - Fetch.AI
Set-Cell "D39" "0.827"
Set-Cell "E39" "  -0.64%  "

# Row 40 - SuiNetwork
Set-Cell "D40" "0.817"
Set-Cell "E40" "  -5.51%  "

# Row 41 - Filecoin
Set-Cell "E41" "  -0.44%  "

# Row 42 - FirstDigitalUSD
Set-Cell "E42" "  +0.06%  "

# Row 43 - Bittensor
Set-Cell "D43" "271.06"
Set-Cell "E43" "  +0.50%  "

# Row 44 - WhiteBITCoin
Set-Cell "E44" "  +0.73%  "

# Row 45 - Stellar
Set-Cell "D45" "0.0954"
Set-Cell "E45" "  +0.60%  "

# Row 46 - Mantle
Set-Cell "D46" "0.589"
Set-Cell "E46" "  +0.43%  "

# Row 47 - Hedera
Set-Cell "E47" "  -0.83%  "

# Row 48 - EnergySwap
Set-Cell "D48" "18.41"
Set-Cell "E48" "  -2.07%  "

# Row 49 - Maker
Set-Cell "D49" "1.966.76"
Set-Cell "E49" "  -0.24%  "

# Row 50 - RenderToken
Set-Cell "D50" "4.52"
Set-Cell "E50" "  -1.84%  "

# Row 51 - VeChain
Set-Cell "E51" "  -0.23%  "
